# Change in unit of AIC: rescale the AIC-related result cells (D7, E5, E7, F8, G7)
# on every year sheet by a factor of 1e-6 (the new run reports the same
# quantities in a different, smaller unit).

$wb = $excel.ActiveWorkbook

$targetCells = @("D7", "E5", "E7", "F8", "G7")
$scale = 0.000001

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($addr in $targetCells) {
        $rng = $ws.Range($addr)
        $val = $rng.Value2
        if ($val -ne 0) {
            $rng.Value2 = $val * $scale
        }
    }
}
